$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N4").Value = 13
$ws.Range("N5").Value = 10
$ws.Range("N6").Value = 11
$ws.Range("N7").Value = 8
$ws.Range("N8").Value = 38

$ws.Range("N4").Select()
